$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated column G ("K") values - replacing former Strike# derived calc
$gValues = @{
    2 = 0
    3 = 1
    4 = 1
    5 = 0
    6 = 0
    7 = 1
    8 = 0
    9 = 0
    10 = 2
    11 = 0
    12 = 2
    13 = 3
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 1
    20 = 0
    21 = 1
    22 = 1
    23 = 1
    24 = 0
    25 = 1
    26 = 0
    27 = 0
    28 = 0
    29 = 0
    30 = 1
    31 = 0
    32 = 1
    33 = 0
    34 = 1
    35 = 1
    36 = 2
    37 = 0
    38 = 0
    39 = 2
    40 = 0
    41 = 0
    42 = 1
    43 = 1
    44 = 0
    45 = 1
    46 = 1
    47 = 0
    48 = 0
    49 = 2
    50 = 0
    51 = 0
    52 = 0
    53 = 1
    54 = 1
    55 = 1
    56 = 1
    57 = 1
    58 = 0
    59 = 2
    60 = 0
    61 = 0
    62 = 1
    63 = 0
    64 = 0
    65 = 1
    66 = 1
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}
